$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) time_variants: insert a new row above row 20 for the ambulatory-care
#    smear-negative intervention ("int_perc_ambulatorycare_smearneg") and
#    push the existing rows 20-37 down to 21-38.
# ---------------------------------------------------------------------------
$wsTV = $wb.Worksheets.Item("time_variants")

$wsTV.Rows.Item(20).Insert()
$wsTV.Rows.Item(20).Clear()

# Re-create the formatting of the new row from existing same-styled cells
# (row 19 carries styles 25/5/12 for A/B/Z..AE, I22 carries style 7).
$wsTV.Range("A19").Copy()
$wsTV.Range("A20").PasteSpecial(-4122)

$wsTV.Range("B19").Copy()
$wsTV.Range("B20").PasteSpecial(-4122)

$wsTV.Range("I22").Copy()
$wsTV.Range("X20").PasteSpecial(-4122)

$wsTV.Range("Z19").Copy()
$wsTV.Range("Z20:AB20").PasteSpecial(-4122)

$wsTV.Range("AD19").Copy()
$wsTV.Range("AD20:AE20").PasteSpecial(-4122)

$wsTV.Range("A20").Value = "int_perc_ambulatorycare_smearneg"
$wsTV.Range("B20").Value = "no"
$wsTV.Range("X20").Value = 0

# ---------------------------------------------------------------------------
# 2) constants: add 5 new rows (53-57) for the ambulatory-care smear-negative
#    intervention cost parameters, matching the style of the rows above
#    (38-47, style 37 on both columns A and B).
# ---------------------------------------------------------------------------
$wsC = $wb.Worksheets.Item("constants")

$wsC.Range("A47:B47").Copy()
$wsC.Range("A53:B57").PasteSpecial(-4122)

$wsC.Range("A53").Value = "econ_unitcost_ambulatorycare_smearneg"
$wsC.Range("B53").Value = 0
$wsC.Range("A54").Value = "econ_inflectioncost_ambulatorycare_smearneg"
$wsC.Range("B54").Value = 0
$wsC.Range("A55").Value = "econ_startupcost_ambulatorycare_smearneg"
$wsC.Range("B55").Value = 0
$wsC.Range("A56").Value = "econ_startupduration_ambulatorycare_smearneg"
$wsC.Range("B56").Value = 0
$wsC.Range("A57").Value = "econ_saturation_ambulatorycare_smearneg"
$wsC.Range("B57").Value = 0

# ---------------------------------------------------------------------------
# 3) View/selection state: the author had been on "constants" (cell G50
#    selected) and switched to "time_variants" (scrolled so column C is the
#    first unfrozen column, cell V21 selected) before saving.
# ---------------------------------------------------------------------------
$wsC.Activate()
$wsC.Range("G50").Select()

$wsTV.Activate()
$wsTV.Range("V21").Select()
